# Breytingar á losunar útreikningum og töflu
$wb = $excel.ActiveWorkbook

$wsAburdur = $wb.Worksheets.Item(1)   # aburðartegundir
$wsLosun   = $wb.Worksheets.Item(3)   # losun

# --- losun sheet: only the selected cell changes ---
$wsLosun.Range("C41").Select()

# --- aburðartegundir sheet ---

# New column header (X1) -> shared string "los_n"
$wsAburdur.Range("X1").Value = "los_n"

# New column X: los_n = H/100*(44/28)*0.01*298, for every data row (2-16)
for ($r = 2; $r -le 16; $r++) {
    $wsAburdur.Range("X$r").Formula = "=H$r/100*(44/28)*0.01*298"
}

# Updated U-column values
$wsAburdur.Range("U2").Value = 1.05
$wsAburdur.Range("U4").Value = 0.0226
$wsAburdur.Range("U5").Value = 0.0226
$wsAburdur.Range("U10").Value = 0.143

# Column widths: A is bestFit-ish 16.71, B:T are 9.14
$wsAburdur.Columns.Item(1).ColumnWidth = 16.7109375
$wsAburdur.Range("B1:T1").EntireColumn.ColumnWidth = 9.140625

# Freeze first column and set the active selection to U3
$wsAburdur.Range("B1").Select()
$excel.ActiveWindow.FreezePanes = $true
$wsAburdur.Range("U3").Select()
